$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped from
# 2023-09-13 (45182) to 2023-09-15 (45184) for every data row (2..360).
$startRow = 2
$endRow = 360
$oldValue = 45182
$newValue = 45184

for ($row = $startRow; $row -le $endRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
